$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row: Fecha(D), Calidad(L), Volumen(M), Precio minimo(N), Precio maximo(O), Precio promedio ponderado(P), Precio $/Kg(S)
$rowData = @{
    2 = @(44403, "Primera", 100, 1200, 1300, 1250, 1250)
    3 = @(44403, "Segunda", 120, 950, 1000, 975, 975)
    4 = @(44326, "Primera", 160, 600, 700, 650, 650)
    5 = @(44379, "Primera", 150, 700, 800, 747, 747)
    6 = @(44379, "Segunda", 140, 500, 600, 543, 543)
    7 = @(44330, "Primera", 200, 1200, 1300, 1250, 1250)
    8 = @(44330, "Segunda", 100, 1000, 1100, 1050, 1050)
    9 = @(44348, "Primera", 120, 1000, 1100, 1050, 1050)
    10 = @(44417, "Primera", 200, 1300, 1400, 1350, 1350)
    11 = @(44407, "Primera", 200, 600, 650, 625, 625)
    12 = @(44309, "Primera", 160, 1400, 1500, 1450, 1450)
    13 = @(44344, "Primera", 140, 1000, 1200, 1100, 1100)
    14 = @(44344, "Segunda", 120, 800, 850, 825, 825)
    15 = @(44350, "Primera", 140, 750, 800, 775, 775)
    16 = @(44414, "Primera", 160, 1300, 1400, 1350, 1350)
    17 = @(44389, "Primera", 140, 750, 800, 775, 775)
    18 = @(44389, "Segunda", 120, 600, 700, 650, 650)
    19 = @(44386, "Primera", 160, 700, 750, 725, 725)
    20 = @(44386, "Segunda", 200, 600, 650, 625, 625)
    21 = @(44316, "Primera", 140, 1100, 1200, 1150, 1150)
    22 = @(44372, "Primera", 900, 750, 800, 772, 772)
    23 = @(44372, "Segunda", 900, 600, 650, 628, 628)
    24 = @(44351, "Primera", 100, 700, 800, 750, 750)
    25 = @(44351, "Segunda", 100, 600, 700, 650, 650)
    26 = @(44358, "Primera", 200, 700, 800, 750, 750)
    27 = @(44358, "Segunda", 200, 600, 650, 625, 625)
    28 = @(44425, "Primera", 140, 1200, 1300, 1250, 1250)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("L$r").Value = $vals[1]
    $ws.Range("M$r").Value = $vals[2]
    $ws.Range("N$r").Value = $vals[3]
    $ws.Range("O$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
    $ws.Range("S$r").Value = $vals[6]
}
